$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1855.25
$ws.Range("I43").Value = 420.5
$ws.Range("J43").Value = 3290
$ws.Range("K43").Value = 420.5
$ws.Range("L43").Value = 3290
$ws.Range("M43").Value = -351.5
$ws.Range("N43").Value = -3428
$ws.Range("H93").Value = 25984.25
$ws.Range("J93").Value = 25984.25
$ws.Range("L93").Value = 25984.25
$ws.Range("N93").Value = -30976.25
$ws.Range("H129").Value = 815.4299999999999
$ws.Range("J129").Value = 866.68134
$ws.Range("L129").Value = 2600.04402
$ws.Range("N129").Value = -12600.04402
$ws.Range("H141").Value = 65134.312
$ws.Range("I141").Value = 79173.38
$ws.Range("K141").Value = 237520.14
$ws.Range("M141").Value = -232340.14

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4591.9014
$ws.Range("I32").Value = 2942.18
$ws.Range("K32").Value = 2942.18
$ws.Range("M32").Value = -2655.18
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41748
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -128736
$ws.Range("H103").Value = 35681
$ws.Range("J103").Value = 35681
$ws.Range("L103").Value = 35681
$ws.Range("N103").Value = -38025
$ws.Range("H122").Value = 4238.222
$ws.Range("I122").Value = 1065.6666
$ws.Range("J122").Value = 10583.333
$ws.Range("K122").Value = 3196.9998
$ws.Range("L122").Value = 31749.999
$ws.Range("M122").Value = -746.9998000000001
$ws.Range("N122").Value = -36649.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5244.75
$ws.Range("I86").Value = 4490
$ws.Range("J86").Value = 5999.5
$ws.Range("K86").Value = 4490
$ws.Range("L86").Value = 5999.5
$ws.Range("M86").Value = -3367
$ws.Range("N86").Value = -8245.5
$ws.Range("H89").Value = 5244.75
$ws.Range("I89").Value = 4490
$ws.Range("J89").Value = 5999.5
$ws.Range("K89").Value = 22450
$ws.Range("L89").Value = 29997.5
$ws.Range("M89").Value = -16834
$ws.Range("N89").Value = -41229.5
$ws.Range("H99").Value = 7695763.5
$ws.Range("I99").Value = 16668380
$ws.Range("J99").Value = 4949.9287
$ws.Range("K99").Value = 16668380
$ws.Range("L99").Value = 4949.9287
$ws.Range("M99").Value = -16666882
$ws.Range("N99").Value = -7945.9287
$ws.Range("H122").Value = 3522.7144
$ws.Range("I122").Value = 1514.75
$ws.Range("J122").Value = 6200
$ws.Range("K122").Value = 4544.25
$ws.Range("L122").Value = 18600
$ws.Range("M122").Value = -2094.25
$ws.Range("N122").Value = -23500
$ws.Range("H126").Value = 7695763.5
$ws.Range("I126").Value = 16668380
$ws.Range("J126").Value = 4949.9287
$ws.Range("K126").Value = 50005140
$ws.Range("L126").Value = 14849.7861
$ws.Range("M126").Value = -50002670
$ws.Range("N126").Value = -19789.7861
$ws.Range("H132").Value = 2328.0637
$ws.Range("I132").Value = 1971.1316
$ws.Range("J132").Value = 3835.111
$ws.Range("K132").Value = 5913.3948
$ws.Range("L132").Value = 11505.333
$ws.Range("M132").Value = -3383.3948
$ws.Range("N132").Value = -16565.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 83334696
$ws.Range("I114").Value = 500000160
$ws.Range("J114").Value = 1609.6
$ws.Range("K114").Value = 1500000480
$ws.Range("L114").Value = 4828.799999999999
$ws.Range("M114").Value = -1499997226
$ws.Range("N114").Value = -11336.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2289.5454
$ws.Range("I102").Value = 1613.48
$ws.Range("K102").Value = 1613.48
$ws.Range("M102").Value = 8.519999999999982
$ws.Range("H107").Value = 11111910
$ws.Range("I107").Value = 540
$ws.Range("J107").Value = 15873926
$ws.Range("K107").Value = 540
$ws.Range("L107").Value = 15873926
$ws.Range("M107").Value = 1380
$ws.Range("N107").Value = -15877766

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4382.5264
$ws.Range("I40").Value = 2697.0908
$ws.Range("J40").Value = 6700
$ws.Range("K40").Value = 2697.0908
$ws.Range("L40").Value = 6700
$ws.Range("M40").Value = -2561.0908
$ws.Range("N40").Value = -6972
$ws.Range("H82").Value = 5792.478
$ws.Range("I82").Value = 7789.0713
$ws.Range("J82").Value = 2686.6667
$ws.Range("K82").Value = 7789.0713
$ws.Range("L82").Value = 2686.6667
$ws.Range("M82").Value = -7428.0713
$ws.Range("N82").Value = -3408.6667
$ws.Range("H85").Value = 5792.478
$ws.Range("I85").Value = 7789.0713
$ws.Range("J85").Value = 2686.6667
$ws.Range("K85").Value = 7789.0713
$ws.Range("L85").Value = 2686.6667
$ws.Range("M85").Value = -6541.0713
$ws.Range("N85").Value = -5182.6667
$ws.Range("H93").Value = 10104636
$ws.Range("I93").Value = 15876285
$ws.Range("J93").Value = 4250
$ws.Range("K93").Value = 15876285
$ws.Range("L93").Value = 4250
$ws.Range("M93").Value = -15875037
$ws.Range("N93").Value = -6746
$ws.Range("H122").Value = 8105.25
$ws.Range("I122").Value = 6159.6665
$ws.Range("J122").Value = 8753.777
$ws.Range("K122").Value = 18478.9995
$ws.Range("L122").Value = 26261.331
$ws.Range("M122").Value = -16028.9995
$ws.Range("N122").Value = -31161.331
$ws.Range("H132").Value = 6504.6665
$ws.Range("I132").Value = 3625.1765
$ws.Range("J132").Value = 11399.8
$ws.Range("K132").Value = 10875.5295
$ws.Range("L132").Value = 34199.39999999999
$ws.Range("M132").Value = -8345.529500000001
$ws.Range("N132").Value = -39259.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 27441.666
$ws.Range("J86").Value = 27441.666
$ws.Range("L86").Value = 27441.666
$ws.Range("N86").Value = -29687.666
$ws.Range("H89").Value = 27441.666
$ws.Range("J89").Value = 27441.666
$ws.Range("L89").Value = 137208.33
$ws.Range("N89").Value = -148440.33
$ws.Range("H107").Value = 664.1
$ws.Range("I107").Value = 649.7778
$ws.Range("J107").Value = 793
$ws.Range("K107").Value = 1949.3334
$ws.Range("L107").Value = 2379
$ws.Range("M107").Value = -29.33339999999998
$ws.Range("N107").Value = -6219
$ws.Range("H113").Value = 7484.0713
$ws.Range("I113").Value = 14454.429
$ws.Range("J113").Value = 513.7143
$ws.Range("K113").Value = 43363.287
$ws.Range("L113").Value = 1541.1429
$ws.Range("M113").Value = -41193.287
$ws.Range("N113").Value = -5881.1429
$ws.Range("H132").Value = 16674086
$ws.Range("I132").Value = 8491.929
$ws.Range("J132").Value = 55560470
$ws.Range("K132").Value = 25475.787
$ws.Range("L132").Value = 166681410
$ws.Range("M132").Value = -22945.787
$ws.Range("N132").Value = -166686470
